# Auto-generated edit script applying the cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.772.63"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "3.392.38"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'570.01"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Value = "'161.28"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.391.57"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("D9").Value = "'0.545"
$ws.Range("E9").Value = "  -5.07%  "
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "'0.419"
$ws.Range("E12").Value = "  -4.81%  "
$ws.Range("D13").Value = "3.977.54"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'26.82"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "63.789.91"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "3.394.53"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").Value = "'13.47"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "'375.68"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'69.93"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").Value = "'0.510"
$ws.Range("E25").Value = "  -4.71%  "
$ws.Range("E26").Value = "  -5.45%  "
$ws.Range("D27").Value = "'9.51"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'6.04"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "'1.39"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'22.71"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").Value = "'7.00"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'159.36"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  -6.25%  "
$ws.Range("E38").Value = "  +8.70%  "
$ws.Range("E39").Value = "  -4.90%  "
$ws.Range("D40").Value = "'0.0718"
$ws.Range("E40").Value = "  -4.10%  "
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").Value = "'42.62"
$ws.Range("D43").Value = "2.734.12"
$ws.Range("E43").Value = "  -5.56%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'6.40"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'25.96"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "'327.36"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("E50").Value = "  -5.05%  "
$ws.Range("E51").Value = "  -2.04%  "
